$wb = $excel.ActiveWorkbook

$newStamp = "February 02 2026 12.49.33 EST"

# --- Sheet "About" ---
$wsAbout = $wb.Worksheets.Item("About")

$wsAbout.Range("A2").Value = "Version: mines - January 30 (built on " + $newStamp + ")"
$wsAbout.Range("A6").Value = "Recommended Citation:  " + [char]34 + "Global Energy Monitor, Coal mine boundaries and methane sources for Pansan Coal Mine, China, M5206, version 'mines - January 30 (built on " + $newStamp + ")'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# --- Sheet "Boundaries and methane sources" ---
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

for ($r = 2; $r -le 13; $r++) {
    $cell = $wsData.Range("S" + $r)
    $cell.Value = "mines - January 30 (built on " + $newStamp + ")"
}
